# Update TPM-derived NATMI ligand-receptor metrics (Il4-Il2rg) with
# newly recomputed values. Columns: E=Ligand-expressing cells,
# F=Ligand detection rate, G=Ligand avg expr, H=Ligand total expr,
# I/J=Ligand derived specificity (avg/total), M=Receptor avg expr,
# N=Receptor total expr, O/P=Receptor derived specificity (avg/total),
# Q/R=Edge expression weight (avg/total), S/T=Edge derived specificity
# (avg/total). Rows 2-10 correspond to the 3x3 ECs/FAPs/MuSCs sending-
# /target-cluster combinations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.3205216666666666
$ws.Cells.Item(2, 8).Value = 0.961565
$ws.Cells.Item(2, 9).Value = 0.1355444593128396
$ws.Cells.Item(2, 10).Value = 0.1355444593128396
$ws.Cells.Item(2, 13).Value = 7.236245333333333
$ws.Cells.Item(2, 14).Value = 21.708736
$ws.Cells.Item(2, 15).Value = 0.6630432242917509
$ws.Cells.Item(2, 16).Value = 0.6630432242917509
$ws.Cells.Item(2, 17).Value = 2.319373414648889
$ws.Cells.Item(2, 18).Value = 20.87436073184
$ws.Cells.Item(2, 19).Value = 0.08987183533766718
$ws.Cells.Item(2, 20).Value = 0.08987183533766718

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.3205216666666666
$ws.Cells.Item(3, 8).Value = 0.961565
$ws.Cells.Item(3, 9).Value = 0.1355444593128396
$ws.Cells.Item(3, 10).Value = 0.1355444593128396
$ws.Cells.Item(3, 15).Value = 0.003787629702975075
$ws.Cells.Item(3, 16).Value = 0.003787629702975075
$ws.Cells.Item(3, 17).Value = 0.013249404135
$ws.Cells.Item(3, 18).Value = 0.119244637215
$ws.Cells.Item(3, 19).Value = 0.0005133922201670077
$ws.Cells.Item(3, 20).Value = 0.0005133922201670077

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.3205216666666666
$ws.Cells.Item(4, 8).Value = 0.961565
$ws.Cells.Item(4, 9).Value = 0.1355444593128396
$ws.Cells.Item(4, 10).Value = 0.1355444593128396
$ws.Cells.Item(4, 13).Value = 3.636103333333333
$ws.Cells.Item(4, 14).Value = 10.90831
$ws.Cells.Item(4, 15).Value = 0.3331691460052741
$ws.Cells.Item(4, 16).Value = 0.3331691460052741
$ws.Cells.Item(4, 17).Value = 1.165449900572222
$ws.Cells.Item(4, 18).Value = 10.48904910515
$ws.Cells.Item(4, 19).Value = 0.04515923175500538
$ws.Cells.Item(4, 20).Value = 0.04515923175500538

$ws.Cells.Item(5, 9).Value = 0.4633580360449179
$ws.Cells.Item(5, 10).Value = 0.4633580360449179
$ws.Cells.Item(5, 13).Value = 7.236245333333333
$ws.Cells.Item(5, 14).Value = 21.708736
$ws.Cells.Item(5, 15).Value = 0.6630432242917509
$ws.Cells.Item(5, 16).Value = 0.6630432242917509
$ws.Cells.Item(5, 17).Value = 7.928766072142221
$ws.Cells.Item(5, 18).Value = 71.35889464927999
$ws.Cells.Item(5, 19).Value = 0.3072264062207157
$ws.Cells.Item(5, 20).Value = 0.3072264062207157

$ws.Cells.Item(6, 9).Value = 0.4633580360449179
$ws.Cells.Item(6, 10).Value = 0.4633580360449179
$ws.Cells.Item(6, 15).Value = 0.003787629702975075
$ws.Cells.Item(6, 16).Value = 0.003787629702975075
$ws.Cells.Item(6, 19).Value = 0.001755028660435927
$ws.Cells.Item(6, 20).Value = 0.001755028660435927

$ws.Cells.Item(7, 9).Value = 0.4633580360449179
$ws.Cells.Item(7, 10).Value = 0.4633580360449179
$ws.Cells.Item(7, 13).Value = 3.636103333333333
$ws.Cells.Item(7, 14).Value = 10.90831
$ws.Cells.Item(7, 15).Value = 0.3331691460052741
$ws.Cells.Item(7, 16).Value = 0.3331691460052741
$ws.Cells.Item(7, 17).Value = 3.984084482505555
$ws.Cells.Item(7, 18).Value = 35.85676034255
$ws.Cells.Item(7, 19).Value = 0.1543766011637663
$ws.Cells.Item(7, 20).Value = 0.1543766011637663

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.9484743333333334
$ws.Cells.Item(8, 8).Value = 2.845423
$ws.Cells.Item(8, 9).Value = 0.4010975046422426
$ws.Cells.Item(8, 10).Value = 0.4010975046422425
$ws.Cells.Item(8, 13).Value = 7.236245333333333
$ws.Cells.Item(8, 14).Value = 21.708736
$ws.Cells.Item(8, 15).Value = 0.6630432242917509
$ws.Cells.Item(8, 16).Value = 0.6630432242917509
$ws.Cells.Item(8, 17).Value = 6.863392968369777
$ws.Cells.Item(8, 18).Value = 61.770536715328
$ws.Cells.Item(8, 19).Value = 0.265944982733368
$ws.Cells.Item(8, 20).Value = 0.265944982733368

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.9484743333333334
$ws.Cells.Item(9, 8).Value = 2.845423
$ws.Cells.Item(9, 9).Value = 0.4010975046422426
$ws.Cells.Item(9, 10).Value = 0.4010975046422425
$ws.Cells.Item(9, 15).Value = 0.003787629702975075
$ws.Cells.Item(9, 16).Value = 0.003787629702975075
$ws.Cells.Item(9, 17).Value = 0.039207083517
$ws.Cells.Item(9, 18).Value = 0.352863751653
$ws.Cells.Item(9, 19).Value = 0.001519208822372141
$ws.Cells.Item(9, 20).Value = 0.001519208822372141

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.9484743333333334
$ws.Cells.Item(10, 8).Value = 2.845423
$ws.Cells.Item(10, 9).Value = 0.4010975046422426
$ws.Cells.Item(10, 10).Value = 0.4010975046422425
$ws.Cells.Item(10, 13).Value = 3.636103333333333
$ws.Cells.Item(10, 14).Value = 10.90831
$ws.Cells.Item(10, 15).Value = 0.3331691460052741
$ws.Cells.Item(10, 16).Value = 0.3331691460052741
$ws.Cells.Item(10, 17).Value = 3.448750685014445
$ws.Cells.Item(10, 18).Value = 31.03875616513
$ws.Cells.Item(10, 19).Value = 0.1336333130865024
$ws.Cells.Item(10, 20).Value = 0.1336333130865024
